# Collections.xlsx example/test data update.
#
# The "DataTable" sample in column D used a single hard-coded "now()"-style
# timestamp (with a date+time number format) for every row. Replace it with
# a short sequence of plain dates (2000-01-01 .. 2000-01-05) formatted as a
# date only, and narrow column D to fit the shorter values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D is narrower now that it only needs to fit a short date.
$ws.Columns("D").ColumnWidth = 9.45

# Display the DataTable date column as a plain date (no time portion).
$ws.Range("D7:D11").NumberFormat = "m/d/yyyy"

# Replace the single repeated timestamp with distinct, deterministic dates.
$ws.Range("D7").Value = 36526   # 1/1/2000
$ws.Range("D8").Value = 36527   # 1/2/2000
$ws.Range("D9").Value = 36528   # 1/3/2000
$ws.Range("D10").Value = 36529  # 1/4/2000
$ws.Range("D11").Value = 36530  # 1/5/2000
